$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meta")

$ws.Range("B1").Value = "2025-09-23 00:24:57"
$ws.Range("B2").Value = "C:\v5\1.pdf"
$ws.Range("A3").Value = "An (din interfață)"

# "2022" looks numeric, so force it to stay text (matches the source
# workbook's inlineStr cell type) and then drop back to the default
# style so no stray number-format gets attached to the cell.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2022"
$ws.Range("B3").Style = "Normal"
